# Update "想去人数" (F column) counts on the 展览 (sheet1) and 全部类型 (sheet4)
# worksheets to reflect refreshed scrape data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Rows are identical between the two sheets for the first 8 updates,
# but diverge afterwards because 全部类型 interleaves rows from other
# sheets, shifting the remaining matching events down by two rows.
$commonUpdates = @(
    @{ Row = 5;  New = 16039 },
    @{ Row = 8;  New = 722 },
    @{ Row = 9;  New = 15526 },
    @{ Row = 11; New = 9149 },
    @{ Row = 14; New = 1023 },
    @{ Row = 20; New = 72 },
    @{ Row = 25; New = 1129 },
    @{ Row = 28; New = 28 }
)

$exhibitOnlyUpdates = @(
    @{ Row = 35; New = 265 },
    @{ Row = 36; New = 337 },
    @{ Row = 39; New = 5628 },
    @{ Row = 40; New = 5236 }
)

$allOnlyUpdates = @(
    @{ Row = 37; New = 265 },
    @{ Row = 38; New = 337 },
    @{ Row = 41; New = 5628 },
    @{ Row = 43; New = 5236 }
)

foreach ($u in $commonUpdates) {
    $wsExhibit.Range("F" + $u.Row).Value = $u.New
    $wsAll.Range("F" + $u.Row).Value = $u.New
}

foreach ($u in $exhibitOnlyUpdates) {
    $wsExhibit.Range("F" + $u.Row).Value = $u.New
}

foreach ($u in $allOnlyUpdates) {
    $wsAll.Range("F" + $u.Row).Value = $u.New
}

$wb.Save()
